$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the three "ECs" sending-cluster rows (old rows 2-4); remaining rows shift up
$ws.Range("A2:A4").EntireRow.Delete()

# Refresh data rows 2-7 with the recomputed TPM-based statistics

# Row 2: FAPs -> ECs
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "F13a1"
$ws.Range("C2").Value = "Itga4"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.07623033333333334
$ws.Range("H2").Value = 0.228691
$ws.Range("I2").Value = 0.7411556909515168
$ws.Range("J2").Value = 0.7411556909515167
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.05968133333333333
$ws.Range("N2").Value = 0.179044
$ws.Range("O2").Value = 0.02602747651633847
$ws.Range("P2").Value = 0.02602747651633848
$ws.Range("Q2").Value = 0.004549527933777778
$ws.Range("R2").Value = 0.040945751404
$ws.Range("S2").Value = 0.01929041234119122
$ws.Range("T2").Value = 0.01929041234119122

# Row 3: FAPs -> FAPs
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "F13a1"
$ws.Range("C3").Value = "Itga4"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.07623033333333334
$ws.Range("H3").Value = 0.228691
$ws.Range("I3").Value = 0.7411556909515168
$ws.Range("J3").Value = 0.7411556909515167
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.3302223333333333
$ws.Range("N3").Value = 0.990667
$ws.Range("O3").Value = 0.144012433133819
$ws.Range("P3").Value = 0.144012433133819
$ws.Range("Q3").Value = 0.02517295854411112
$ws.Range("R3").Value = 0.226556626897
$ws.Range("S3").Value = 0.1067356343849047
$ws.Range("T3").Value = 0.1067356343849047

# Row 4: FAPs -> MuSCs
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "F13a1"
$ws.Range("C4").Value = "Itga4"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.07623033333333334
$ws.Range("H4").Value = 0.228691
$ws.Range("I4").Value = 0.7411556909515168
$ws.Range("J4").Value = 0.7411556909515167
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.903109
$ws.Range("N4").Value = 5.709327
$ws.Range("O4").Value = 0.8299600903498424
$ws.Range("P4").Value = 0.8299600903498425
$ws.Range("Q4").Value = 0.1450746334396667
$ws.Range("R4").Value = 1.305671700957
$ws.Range("S4").Value = 0.6151296442254207
$ws.Range("T4").Value = 0.6151296442254207

# Row 5: MuSCs -> ECs
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "F13a1"
$ws.Range("C5").Value = "Itga4"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.026623
$ws.Range("H5").Value = 0.079869
$ws.Range("I5").Value = 0.2588443090484832
$ws.Range("J5").Value = 0.2588443090484832
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.05968133333333333
$ws.Range("N5").Value = 0.179044
$ws.Range("O5").Value = 0.02602747651633847
$ws.Range("P5").Value = 0.02602747651633848
$ws.Range("Q5").Value = 0.001588896137333333
$ws.Range("R5").Value = 0.014300065236
$ws.Range("S5").Value = 0.006737064175147255
$ws.Range("T5").Value = 0.006737064175147256

# Row 6: MuSCs -> FAPs
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "F13a1"
$ws.Range("C6").Value = "Itga4"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.026623
$ws.Range("H6").Value = 0.079869
$ws.Range("I6").Value = 0.2588443090484832
$ws.Range("J6").Value = 0.2588443090484832
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.3302223333333333
$ws.Range("N6").Value = 0.990667
$ws.Range("O6").Value = 0.144012433133819
$ws.Range("P6").Value = 0.144012433133819
$ws.Range("Q6").Value = 0.008791509180333333
$ws.Range("R6").Value = 0.079123582623
$ws.Range("S6").Value = 0.03727679874891427
$ws.Range("T6").Value = 0.03727679874891427

# Row 7: MuSCs -> MuSCs
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "F13a1"
$ws.Range("C7").Value = "Itga4"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.026623
$ws.Range("H7").Value = 0.079869
$ws.Range("I7").Value = 0.2588443090484832
$ws.Range("J7").Value = 0.2588443090484832
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.903109
$ws.Range("N7").Value = 5.709327
$ws.Range("O7").Value = 0.8299600903498424
$ws.Range("P7").Value = 0.8299600903498425
$ws.Range("Q7").Value = 0.05066647090699999
$ws.Range("R7").Value = 0.455998238163
$ws.Range("S7").Value = 0.2148304461244216
$ws.Range("T7").Value = 0.2148304461244217
